# Bump the embedded "term" metadata to the 1.1.0 release: update the
# Version and Date rows on the Metadata sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "1.1.0"
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
